# Apply updated F-column ('想去人数' / interest count) values across sheets
# Sheet index -> name mapping: 1=展览, 2=演出, 3=本地生活, 4=全部类型
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 76
$ws.Range("F3").Value = 165
$ws.Range("F5").Value = 50
$ws.Range("F6").Value = 2744
$ws.Range("F8").Value = 1625
$ws.Range("F9").Value = 7431
$ws.Range("F11").Value = 7620
$ws.Range("F15").Value = 6128
$ws.Range("F16").Value = 3248
$ws.Range("F17").Value = 3622
$ws.Range("F18").Value = 14
$ws.Range("F19").Value = 9
$ws.Range("F20").Value = 15
$ws.Range("F22").Value = 443
$ws.Range("F24").Value = 281
$ws.Range("F25").Value = 281
$ws.Range("F26").Value = 3613
$ws.Range("F27").Value = 113
$ws.Range("F28").Value = 337
$ws.Range("F29").Value = 925
$ws.Range("F30").Value = 256
$ws.Range("F31").Value = 1080
$ws.Range("F32").Value = 61
$ws.Range("F33").Value = 16
$ws.Range("F34").Value = 2601
$ws.Range("F35").Value = 1450
$ws.Range("F36").Value = 10
$ws.Range("F37").Value = 13
$ws.Range("F38").Value = 21
$ws.Range("F39").Value = 3235
$ws.Range("F40").Value = 151
$ws.Range("F41").Value = 240
$ws.Range("F44").Value = 476
$ws.Range("F45").Value = 1267
$ws.Range("F48").Value = 585

$ws = $wb.Worksheets.Item(2)
$ws.Range("F5").Value = 233
$ws.Range("F8").Value = 38
$ws.Range("F9").Value = 398

$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 119

$ws = $wb.Worksheets.Item(4)
$ws.Range("F4").Value = 76
$ws.Range("F5").Value = 165
$ws.Range("F7").Value = 50
$ws.Range("F8").Value = 119
$ws.Range("F9").Value = 2744
$ws.Range("F10").Value = 1625
$ws.Range("F11").Value = 233
$ws.Range("F13").Value = 7431
$ws.Range("F14").Value = 7620
$ws.Range("F17").Value = 6128
$ws.Range("F18").Value = 3248
$ws.Range("F19").Value = 3622
$ws.Range("F20").Value = 14
$ws.Range("F21").Value = 9
$ws.Range("F23").Value = 443
$ws.Range("F26").Value = 281
$ws.Range("F28").Value = 281
$ws.Range("F29").Value = 3613
$ws.Range("F31").Value = 113
$ws.Range("F34").Value = 337
$ws.Range("F35").Value = 925
$ws.Range("F36").Value = 256
$ws.Range("F37").Value = 16
$ws.Range("F38").Value = 2601
$ws.Range("F39").Value = 1450
$ws.Range("F40").Value = 10
$ws.Range("F41").Value = 13
$ws.Range("F43").Value = 3235
$ws.Range("F44").Value = 240
$ws.Range("F46").Value = 476
$ws.Range("F47").Value = 1267
$ws.Range("F49").Value = 585
